$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ficha técnica")

# Original "Ficha técnica" rows (A:B):
#   1: "" / " "
#   2: DERECHO / Vivienda
#   3: DIMENSIÓN / Habitabilidad        <- remove this row
#   4: CONINDICADOR / Hacinamiento
#   5: NOMINDICADOR / Porcentaje...
#   6: DEFINICIÓN / El indicador mide...
#   7: CÁLCULO / Para cada año calcular...
#
# New layout: drop the DIMENSIÓN row (rows below shift up one), then append
# two new rows (TIPOIND / Resultados and CITA / UMAD con base en...).

$ws.Rows.Item(3).Delete()

$ws.Range("A7").Value = "TIPOIND"
$ws.Range("B7").Value = "Resultados"
$ws.Range("A8").Value = "CITA"
$ws.Range("B8").Value = "UMAD con base en Instituto de Economía, Universidad de la República (2020) Encuesta Continua de Hogares Compatibilizada 1981-2018 Versión 12 DOI: http://doiorg/1047426/ECHINE"
